$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay text (preserve trailing zeros / formatting)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "20.153.05"
$ws.Cells.Item(2, 5).Value = "  +1.02%  "

$ws.Cells.Item(3, 4).Value = "1.432.48"
$ws.Cells.Item(3, 5).Value = "  +1.46%  "

$ws.Cells.Item(4, 4).Value = "1.019"
$ws.Cells.Item(4, 5).Value = "  +1.79%  "

$ws.Cells.Item(5, 4).Value = "276.60"
$ws.Cells.Item(5, 5).Value = "  +1.15%  "

$ws.Cells.Item(6, 4).Value = "0.8981"
$ws.Cells.Item(6, 5).Value = "  -10.27%  "

$ws.Cells.Item(7, 4).Value = "0.3671"
$ws.Cells.Item(7, 5).Value = "  -0.96%  "

$ws.Cells.Item(8, 4).Value = "0.3127"
$ws.Cells.Item(8, 5).Value = "  +1.81%  "

$ws.Cells.Item(9, 4).Value = "38.95"
$ws.Cells.Item(9, 5).Value = "  -0.79%  "

$ws.Cells.Item(10, 4).Value = "1.010"
$ws.Cells.Item(10, 5).Value = "  +1.28%  "

$ws.Cells.Item(11, 4).Value = "0.06501"
$ws.Cells.Item(11, 5).Value = "  -0.94%  "

$ws.Cells.Item(12, 4).Value = "1.013"
$ws.Cells.Item(12, 5).Value = "  +1.13%  "

$ws.Cells.Item(13, 4).Value = "5.392"
$ws.Cells.Item(13, 5).Value = "  +0.39%  "

$ws.Cells.Item(14, 4).Value = "17.25"
$ws.Cells.Item(14, 5).Value = "  +1.61%  "

$ws.Cells.Item(15, 4).Value = "6.072"
$ws.Cells.Item(15, 5).Value = "  -1.67%  "

$ws.Cells.Item(16, 4).Value = "1.445.92"
$ws.Cells.Item(16, 5).Value = "  +2.44%  "

$ws.Cells.Item(17, 4).Value = "0.00001015"
$ws.Cells.Item(17, 5).Value = "  +0.74%  "

$ws.Cells.Item(18, 4).Value = "0.05608"
$ws.Cells.Item(18, 5).Value = "  -2.87%  "

$ws.Cells.Item(19, 4).Value = "0.9050"
$ws.Cells.Item(19, 5).Value = "  -9.56%  "

$ws.Cells.Item(20, 4).Value = "66.71"
$ws.Cells.Item(20, 5).Value = "  -9.34%  "

$ws.Cells.Item(21, 4).Value = "5.449"
$ws.Cells.Item(21, 5).Value = "  -2.81%  "

$ws.Cells.Item(22, 4).Value = "14.32"
$ws.Cells.Item(22, 5).Value = "  -0.71%  "

$ws.Cells.Item(23, 4).Value = "10.93"
$ws.Cells.Item(23, 5).Value = "  +0.52%  "

$ws.Cells.Item(24, 4).Value = "2.263"
$ws.Cells.Item(24, 5).Value = "  -0.75%  "

$ws.Cells.Item(25, 4).Value = "20.312.31"
$ws.Cells.Item(25, 5).Value = "  +1.73%  "

$ws.Cells.Item(26, 4).Value = "2.177"
$ws.Cells.Item(26, 5).Value = "  -3.63%  "

$ws.Cells.Item(27, 4).Value = "135.42"
$ws.Cells.Item(27, 5).Value = "  -2.60%  "

$ws.Cells.Item(28, 4).Value = "16.93"
$ws.Cells.Item(28, 5).Value = "  +0.59%  "

$ws.Cells.Item(29, 4).Value = "1.601.60"
$ws.Cells.Item(29, 5).Value = "  +2.06%  "

$ws.Cells.Item(30, 4).Value = "110.82"
$ws.Cells.Item(30, 5).Value = "  +1.50%  "

$ws.Cells.Item(31, 4).Value = "3.622"
$ws.Cells.Item(31, 5).Value = "  -4.71%  "

$ws.Cells.Item(32, 2).Value = "Filecoin"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(32, 4).Value = "4.869"
$ws.Cells.Item(32, 5).Value = "  -9.08%  "

$ws.Cells.Item(33, 2).Value = "ImmutableX"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(33, 4).Value = "0.7959"
$ws.Cells.Item(33, 5).Value = "  -4.02%  "

$ws.Cells.Item(34, 4).Value = "0.07644"
$ws.Cells.Item(34, 5).Value = "  -0.61%  "

$ws.Cells.Item(35, 4).Value = "0.05913"
$ws.Cells.Item(35, 5).Value = "  +2.35%  "

$ws.Cells.Item(36, 4).Value = "1.425"
$ws.Cells.Item(36, 5).Value = "  +11.70%  "

$ws.Cells.Item(37, 4).Value = "4.686"
$ws.Cells.Item(37, 5).Value = "  -2.52%  "

$ws.Cells.Item(38, 4).Value = "1.114"
$ws.Cells.Item(38, 5).Value = "  +4.72%  "

$ws.Cells.Item(39, 4).Value = "0.01996"
$ws.Cells.Item(39, 5).Value = "  -2.40%  "

$ws.Cells.Item(40, 4).Value = "10.24"
$ws.Cells.Item(40, 5).Value = "  +0.89%  "

$ws.Cells.Item(41, 4).Value = "0.1824"
$ws.Cells.Item(41, 5).Value = "  -5.64%  "

$ws.Cells.Item(42, 4).Value = "0.9134"
$ws.Cells.Item(42, 5).Value = "  -8.70%  "

$ws.Cells.Item(43, 4).Value = "3.533"
$ws.Cells.Item(43, 5).Value = "  +0.09%  "

$ws.Cells.Item(44, 4).Value = "0.5231"
$ws.Cells.Item(44, 5).Value = "  -1.23%  "

$ws.Cells.Item(45, 2).Value = "EnergySwap"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(45, 4).Value = "12.03"
$ws.Cells.Item(45, 5).Value = "  -1.57%  "

$ws.Cells.Item(46, 2).Value = "FraxShare"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(46, 4).Value = "6.719"
$ws.Cells.Item(46, 5).Value = "  -20.21%  "

$ws.Cells.Item(47, 4).Value = "119.33"
$ws.Cells.Item(47, 5).Value = "  +8.01%  "

$ws.Cells.Item(48, 4).Value = "0.5119"
$ws.Cells.Item(48, 5).Value = "  -0.05%  "

$ws.Cells.Item(49, 4).Value = "1.754"
$ws.Cells.Item(49, 5).Value = "  -2.81%  "

$ws.Cells.Item(50, 4).Value = "0.06297"
$ws.Cells.Item(50, 5).Value = "  +2.06%  "

$ws.Cells.Item(51, 4).Value = "0.9994"
$ws.Cells.Item(51, 5).Value = "  -0.13%  "
